# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 (swap of rows 42/43 content, keeping column A as row index) ---
$ws.Range("B42").Value = 7165060
$ws.Range("F42").Value = "Arouca"
$ws.Range("G42").Value = "Casa Pia"
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = "A"
$ws.Range("K42").Value = 2.25
$ws.Range("L42").Value = 3.25
$ws.Range("N42").Value = 2.8
$ws.Range("O42").Value = 3.2
$ws.Range("P42").Value = 2.55
$ws.Range("Q42").Value = 0
$ws.Range("T42").Value = 2.25
$ws.Range("U42").Value = 1.975
$ws.Range("V42").Value = 1.875
$ws.Range("W42").Value = -1
$ws.Range("Y42").Value = 1.55
$ws.Range("Z42").Value = -1
$ws.Range("AA42").Value = 0.8500000000000001
$ws.Range("AB42").Value = -1
$ws.Range("AC42").Value = 0.875

# --- Row 43 ---
$ws.Range("B43").Value = 7165059
$ws.Range("F43").Value = "Gil Vicente"
$ws.Range("G43").Value = "Estoril"
$ws.Range("H43").Value = 5
$ws.Range("I43").Value = 3
$ws.Range("J43").Value = "H"
$ws.Range("K43").Value = 2.15
$ws.Range("L43").Value = 3.4
$ws.Range("N43").Value = 2
$ws.Range("O43").Value = 3.5
$ws.Range("P43").Value = 3.6
$ws.Range("Q43").Value = -0.5
$ws.Range("T43").Value = 2.5
$ws.Range("U43").Value = 1.925
$ws.Range("V43").Value = 1.925
$ws.Range("W43").Value = 1
$ws.Range("Y43").Value = -1
$ws.Range("Z43").Value = 1
$ws.Range("AA43").Value = -1
$ws.Range("AB43").Value = 0.925
$ws.Range("AC43").Value = -1

# --- Row 151 (swap of rows 151/152 content) ---
$ws.Range("B151").Value = 6876586
$ws.Range("F151").Value = "Benfica"
$ws.Range("G151").Value = "Rio Ave"
$ws.Range("H151").Value = 4
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = "H"
$ws.Range("K151").Value = 1.166
$ws.Range("L151").Value = 7.5
$ws.Range("M151").Value = 15
$ws.Range("N151").Value = 1.2
$ws.Range("O151").Value = 8
$ws.Range("P151").Value = 10
$ws.Range("Q151").Value = -2
$ws.Range("T151").Value = 3.25
$ws.Range("U151").Value = 1.925
$ws.Range("V151").Value = 1.925
$ws.Range("W151").Value = 0.2
$ws.Range("Y151").Value = -1
$ws.Range("Z151").Value = 1.025
$ws.Range("AA151").Value = -1
$ws.Range("AB151").Value = 0.925

# --- Row 152 ---
$ws.Range("B152").Value = 6876591
$ws.Range("F152").Value = "Vizela"
$ws.Range("G152").Value = "Boavista"
$ws.Range("H152").Value = 1
$ws.Range("I152").Value = 4
$ws.Range("J152").Value = "A"
$ws.Range("K152").Value = 2.3
$ws.Range("L152").Value = 3.1
$ws.Range("M152").Value = 3.25
$ws.Range("N152").Value = 1.95
$ws.Range("O152").Value = 3.2
$ws.Range("P152").Value = 4.2
$ws.Range("Q152").Value = -0.5
$ws.Range("T152").Value = 2.25
$ws.Range("U152").Value = 2
$ws.Range("V152").Value = 1.85
$ws.Range("W152").Value = -1
$ws.Range("Y152").Value = 3.2
$ws.Range("Z152").Value = -1
$ws.Range("AA152").Value = 0.825
$ws.Range("AB152").Value = 1

# --- Row 167 (swap of rows 167/168 content) ---
$ws.Range("B167").Value = 6876603
$ws.Range("F167").Value = "Gil Vicente"
$ws.Range("G167").Value = "Guimaraes"
$ws.Range("H167").Value = 1
$ws.Range("K167").Value = 3.1
$ws.Range("M167").Value = 2.25
$ws.Range("N167").Value = 3.25
$ws.Range("O167").Value = 3.3
$ws.Range("P167").Value = 2.25
$ws.Range("Q167").Value = 0.25
$ws.Range("R167").Value = 1.925
$ws.Range("S167").Value = 1.925
$ws.Range("T167").Value = 2.25
$ws.Range("U167").Value = 1.925
$ws.Range("V167").Value = 1.925
$ws.Range("W167").Value = 2.25
$ws.Range("Z167").Value = 0.925
$ws.Range("AB167").Value = -1
$ws.Range("AC167").Value = 0.925

# --- Row 168 ---
$ws.Range("B168").Value = 6876607
$ws.Range("F168").Value = "Arouca"
$ws.Range("G168").Value = "Vizela"
$ws.Range("H168").Value = 5
$ws.Range("K168").Value = 2
$ws.Range("M168").Value = 3.75
$ws.Range("N168").Value = 1.8
$ws.Range("O168").Value = 3.8
$ws.Range("P168").Value = 4.2
$ws.Range("Q168").Value = -0.75
$ws.Range("R168").Value = 2.05
$ws.Range("S168").Value = 1.8
$ws.Range("T168").Value = 2.5
$ws.Range("U168").Value = 1.85
$ws.Range("V168").Value = 2
$ws.Range("W168").Value = 0.8
$ws.Range("Z168").Value = 1.05
$ws.Range("AB168").Value = 0.8500000000000001
$ws.Range("AC168").Value = -1

# --- Rows 253 / 255-261: individual odds updates ---
$ws.Range("O253").Value = 6

$ws.Range("N255").Value = 1.25
$ws.Range("O255").Value = 6.5
$ws.Range("Q255").Value = -1.75
$ws.Range("R255").Value = 2.03
$ws.Range("S255").Value = 1.87
$ws.Range("U255").Value = 1.875
$ws.Range("V255").Value = 1.975

$ws.Range("N256").Value = 4.5
$ws.Range("P256").Value = 1.727

$ws.Range("N257").Value = 2.5
$ws.Range("P257").Value = 2.9

$ws.Range("O258").Value = 4
$ws.Range("P258").Value = 5.5
$ws.Range("R258").Value = 1.82
$ws.Range("S258").Value = 2.08

$ws.Range("N259").Value = 2.7
$ws.Range("O259").Value = 3.25
$ws.Range("P259").Value = 2.7
$ws.Range("U259").Value = 2
$ws.Range("V259").Value = 1.85

$ws.Range("N260").Value = 1.285
$ws.Range("P260").Value = 11

$ws.Range("O261").Value = 3.5
$ws.Range("P261").Value = 4.2
$ws.Range("R261").Value = 1.93
$ws.Range("S261").Value = 1.97
$ws.Range("U261").Value = 1.95
$ws.Range("V261").Value = 1.9
